$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update contrast labels (column B) for rows 2-11
$ws.Range("B2").Value  = "(Deep/Low SAV) / (Exposed/Low SAV)"
$ws.Range("B3").Value  = "(Exposed/Low SAV) / (Mod/Dense SAV)"
$ws.Range("B4").Value  = "(Exposed/Low SAV) / (Shallow/Dense SAV)"
$ws.Range("B5").Value  = "(Deep/Low SAV) / (Shallow/Low SAV)"
$ws.Range("B6").Value  = "(Exposed/Low SAV) / (Shallow/Low SAV)"
$ws.Range("B7").Value  = "(Mod/Dense SAV) / (Shallow/Low SAV)"
$ws.Range("B8").Value  = "(Shallow/Dense SAV) / (Shallow/Low SAV)"
$ws.Range("B9").Value  = "(Deep/Low SAV) / (Mod/Dense SAV)"
$ws.Range("B10").Value = "(Deep/Low SAV) / (Shallow/Dense SAV)"
$ws.Range("B11").Value = "(Mod/Dense SAV) / (Shallow/Dense SAV)"

# Update ratio (D), std_error (E), statistic (H), adj_p_value (I) for rows 2-11
$ws.Range("D2").Value = 0.792764447167322
$ws.Range("E2").Value = 0.0122161963586223
$ws.Range("H2").Value = -15.0704033942816
$ws.Range("I2").Value = 0.0000000000000000000000000000000000000000000000000253560129516999

$ws.Range("D3").Value = 1.2442487362932
$ws.Range("E3").Value = 0.0198818323656434
$ws.Range("H3").Value = 13.6762077146346
$ws.Range("I3").Value = 0.0000000000000000000000000000000000000000140849157127503

$ws.Range("D4").Value = 1.25959786863902
$ws.Range("E4").Value = 0.0216322800537485
$ws.Range("H4").Value = 13.4385170327407
$ws.Range("I4").Value = 0.000000000000000000000000000000000000000359570400077847

$ws.Range("D5").Value = 0.947952594928598
$ws.Range("E5").Value = 0.0039164362623695
$ws.Range("H5").Value = -12.9374781973715
$ws.Range("I5").Value = 0.00000000000000000000000000000000000027655926880218

$ws.Range("D6").Value = 1.19575568545712
$ws.Range("E6").Value = 0.0187025858357287
$ws.Range("H6").Value = 11.4302503534364
$ws.Range("I6").Value = 0.0000000000000000000000000000295256326507579

$ws.Range("D7").Value = 0.961026240636946
$ws.Range("E7").Value = 0.00631926663580832
$ws.Range("H7").Value = -6.04567288632263
$ws.Range("I7").Value = 0.0000000148787703422299

$ws.Range("D8").Value = 0.949315424572066
$ws.Range("E8").Value = 0.00843461780459183
$ws.Range("H8").Value = -5.8541887015342
$ws.Range("I8").Value = 0.0000000479344132686937

$ws.Range("D9").Value = 0.986396161566116
$ws.Range("E9").Value = 0.00539156003800784
$ws.Range("H9").Value = -2.50593217051231
$ws.Range("I9").Value = 0.122129046086697

$ws.Range("D10").Value = 0.998564407984751
$ws.Range("E10").Value = 0.00794736347564217
$ws.Range("H10").Value = -0.180507795316859
$ws.Range("I10").Value = 1

$ws.Range("D11").Value = 1.01233606424351
$ws.Range("E11").Value = 0.00862212494659611
$ws.Range("H11").Value = 1.43953405871573
$ws.Range("I11").Value = 1
